$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6750
$ws.Range("I62").Value = 5666.6665
$ws.Range("K62").Value = 5666.6665
$ws.Range("M62").Value = -5042.6665
$ws.Range("H65").Value = 6750
$ws.Range("I65").Value = 5666.6665
$ws.Range("K65").Value = 28333.3325
$ws.Range("M65").Value = -25213.3325
$ws.Range("H92").Value = 1116.7916
$ws.Range("I92").Value = 1115.9474
$ws.Range("K92").Value = 1115.9474
$ws.Range("M92").Value = 132.0526
$ws.Range("H100").Value = 1648.9565
$ws.Range("I100").Value = 1221.875
$ws.Range("J100").Value = 1876.7333
$ws.Range("K100").Value = 1221.875
$ws.Range("L100").Value = 1876.7333
$ws.Range("M100").Value = -680.875
$ws.Range("N100").Value = -2958.7333
$ws.Range("H106").Value = 12265.333
$ws.Range("I106").Value = 13867.8
$ws.Range("K106").Value = 13867.8
$ws.Range("M106").Value = -13236.8
$ws.Range("H132").Value = 7579670.5
$ws.Range("I132").Value = 9263411
$ws.Range("J132").Value = 2839
$ws.Range("K132").Value = 27790233
$ws.Range("L132").Value = 8517
$ws.Range("M132").Value = -27787703
$ws.Range("N132").Value = -13577
$ws.Range("H140").Value = 34051.43
$ws.Range("J140").Value = 34051.43
$ws.Range("L140").Value = 34051.43
$ws.Range("N140").Value = -44411.43
$ws.Range("H141").Value = 580.82355
$ws.Range("I141").Value = 580.82355
$ws.Range("K141").Value = 1742.47065
$ws.Range("M141").Value = 3437.52935

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3166.7356
$ws.Range("I32").Value = 2862.2375
$ws.Range("K32").Value = 2862.2375
$ws.Range("M32").Value = -2575.2375
$ws.Range("H35").Value = 2321
$ws.Range("I35").Value = 2321
$ws.Range("K35").Value = 2321
$ws.Range("M35").Value = -1915
$ws.Range("H61").Value = 1921.1111
$ws.Range("I61").Value = 1755.7142
$ws.Range("K61").Value = 1755.7142
$ws.Range("M61").Value = -1543.7142
$ws.Range("H132").Value = 1800.9166
$ws.Range("I132").Value = 1753.8334
$ws.Range("J132").Value = 1848
$ws.Range("K132").Value = 5261.5002
$ws.Range("L132").Value = 5544
$ws.Range("M132").Value = -2731.5002
$ws.Range("N132").Value = -10604
$ws.Range("H136").Value = 1921.1111
$ws.Range("I136").Value = 1755.7142
$ws.Range("K136").Value = 5267.142599999999
$ws.Range("M136").Value = -2717.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1756.6842
$ws.Range("J20").Value = 1703.375
$ws.Range("L20").Value = 1703.375
$ws.Range("N20").Value = -2197.375
$ws.Range("H99").Value = 22728058
$ws.Range("I99").Value = 26316472
$ws.Range("K99").Value = 26316472
$ws.Range("M99").Value = -26314974
$ws.Range("H105").Value = 100990184
$ws.Range("I105").Value = 100990184
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 100990184
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -100988437
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1497.091
$ws.Range("I31").Value = 1218.6666
$ws.Range("J31").Value = 2750
$ws.Range("K31").Value = 1218.6666
$ws.Range("L31").Value = 2750
$ws.Range("M31").Value = -923.6666
$ws.Range("N31").Value = -3340
$ws.Range("H34").Value = 1497.091
$ws.Range("I34").Value = 1218.6666
$ws.Range("J34").Value = 2750
$ws.Range("K34").Value = 1218.6666
$ws.Range("L34").Value = 2750
$ws.Range("M34").Value = -1016.6666
$ws.Range("N34").Value = -3154
$ws.Range("H58").Value = 1469.3334
$ws.Range("I58").Value = 1469.3334
$ws.Range("K58").Value = 1469.3334
$ws.Range("M58").Value = -1266.3334
$ws.Range("H132").Value = 1624.4482
$ws.Range("I132").Value = 914.13635
$ws.Range("K132").Value = 2742.40905
$ws.Range("M132").Value = -212.4090500000002
$ws.Range("H134").Value = 857
$ws.Range("I134").Value = 746.75
$ws.Range("J134").Value = 1209.8
$ws.Range("K134").Value = 2240.25
$ws.Range("L134").Value = 3629.4
$ws.Range("M134").Value = 294.75
$ws.Range("N134").Value = -8699.4
$ws.Range("H136").Value = 1469.3334
$ws.Range("I136").Value = 1469.3334
$ws.Range("K136").Value = 4408.0002
$ws.Range("M136").Value = -1858.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1388.9678
$ws.Range("J5").Value = 767.8
$ws.Range("L5").Value = 2303.4
$ws.Range("N5").Value = -2527.4
$ws.Range("H32").Value = 1966.8889
$ws.Range("J32").Value = 2125
$ws.Range("L32").Value = 6375
$ws.Range("N32").Value = -6941
$ws.Range("H36").Value = 400
$ws.Range("I36").Value = 400
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1200
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1031
$ws.Range("N36").ClearContents()
$ws.Range("H131").Value = 10001972
$ws.Range("I131").Value = 142857400
$ws.Range("J131").Value = 2100.14
$ws.Range("K131").Value = 428572200
$ws.Range("L131").Value = 6300.42
$ws.Range("M131").Value = -428567160
$ws.Range("N131").Value = -16380.42
$ws.Range("H133").Value = 2444.3333
$ws.Range("I133").Value = 1700
$ws.Range("K133").Value = 5100
$ws.Range("M133").Value = -40
$ws.Range("H135").Value = 1388.9678
$ws.Range("J135").Value = 767.8
$ws.Range("L135").Value = 6910.2
$ws.Range("N135").Value = -11980.2
$ws.Range("H141").Value = 4250

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2398.9565
$ws.Range("I132").Value = 1537.5385
$ws.Range("J132").Value = 3518.8
$ws.Range("K132").Value = 4612.6155
$ws.Range("L132").Value = 10556.4
$ws.Range("M132").Value = -2082.6155
$ws.Range("N132").Value = -15616.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 30197
$ws.Range("J96").Value = 30197
$ws.Range("L96").Value = 30197
$ws.Range("N96").Value = -35689
$ws.Range("H104").Value = 4997.375
$ws.Range("J104").Value = 4997.375
$ws.Range("L104").Value = 4997.375
$ws.Range("N104").Value = -11985.375
$ws.Range("H122").Value = 13891133
$ws.Range("I122").Value = 22728946
$ws.Range("J122").Value = 3142.8572
$ws.Range("K122").Value = 68186838
$ws.Range("L122").Value = 9428.571599999999
$ws.Range("M122").Value = -68184388
$ws.Range("N122").Value = -14328.5716
$ws.Range("H132").Value = 25571.953
$ws.Range("I132").Value = 1110
$ws.Range("J132").Value = 52480.1
$ws.Range("K132").Value = 3330
$ws.Range("L132").Value = 157440.3
$ws.Range("M132").Value = -800
$ws.Range("N132").Value = -162500.3
$ws.Range("H136").Value = 2193.0908
$ws.Range("I136").Value = 2235.4443
$ws.Range("J136").Value = 2002.5
$ws.Range("K136").Value = 6706.3329
$ws.Range("L136").Value = 6007.5
$ws.Range("M136").Value = -4156.3329
$ws.Range("N136").Value = -11107.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 30330
$ws.Range("J94").Value = 30330
$ws.Range("L94").Value = 30330
$ws.Range("N94").Value = -32132
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H126").Value = 43479450
$ws.Range("I126").Value = 50001224
$ws.Range("J126").Value = 930
$ws.Range("K126").Value = 150003672
$ws.Range("L126").Value = 2790
$ws.Range("M126").Value = -150001202
$ws.Range("N126").Value = -7730
$ws.Range("H132").Value = 1735.2727
$ws.Range("I132").Value = 1383.0526
$ws.Range("J132").Value = 3966
$ws.Range("K132").Value = 4149.1578
$ws.Range("L132").Value = 11898
$ws.Range("M132").Value = -1619.1578
$ws.Range("N132").Value = -16958
